# Chapter13-HW-3-withAnswers.xlsx edit
#
# 1. Rename the first worksheet "2-hiddenLayers" -> "NN".
#    (Excel automatically rewrites every defined-name formula that
#    referenced '2-hiddenLayers'!... to the new, unquoted NN!... form.)
# 2. Move the active selection on that sheet from P3 to W6.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "NN"

$ws.Activate()
$ws.Range("W6").Select()
